# Swap the order of "Recorded By" names in column G:
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#
# This mirrors a commit that reordered the recorder names in the
# attendance report's "Recorded By" column (col G) on the
# "Session Analysis Results" sheet, for every row whose value matches
# exactly "System, dnasr281@gmail.com" (rows that only contain "System"
# or only "dnasr281@gmail.com" are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

# Find the last used row so we cover the whole table.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value()
    if ($val -eq $oldText) {
        $cell.Value = $newText
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cells in column G"
